$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 328, shifting existing rows 328:390 down to 331:393
$ws.Rows("328:330").Insert()

# Constant columns shared by every data row in this sheet
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100102
$producto = "Cítricos"
$categoriaId = 100102004
$categoria = "Mandarina"

# New row 328: Clementina / Especial
$ws.Cells.Item(328, 1).Value = $mercadoId
$ws.Cells.Item(328, 2).Value = $mercado
$ws.Cells.Item(328, 3).Value = $region
$ws.Cells.Item(328, 4).Value = 44504
$ws.Cells.Item(328, 5).Value = $codreg
$ws.Cells.Item(328, 6).Value = $tipo
$ws.Cells.Item(328, 7).Value = $productoId
$ws.Cells.Item(328, 8).Value = $producto
$ws.Cells.Item(328, 9).Value = $categoriaId
$ws.Cells.Item(328, 10).Value = $categoria
$ws.Cells.Item(328, 11).Value = "Clementina"
$ws.Cells.Item(328, 12).Value = "Especial"
$ws.Cells.Item(328, 13).Value = 350
$ws.Cells.Item(328, 14).Value = 6000
$ws.Cells.Item(328, 15).Value = 6000
$ws.Cells.Item(328, 16).Value = 6000
$ws.Cells.Item(328, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(328, 18).Value = "Provincia de Petorca"
$ws.Cells.Item(328, 19).Value = 600
$ws.Cells.Item(328, 20).Value = 10

# New row 329: Clementina / Primera
$ws.Cells.Item(329, 1).Value = $mercadoId
$ws.Cells.Item(329, 2).Value = $mercado
$ws.Cells.Item(329, 3).Value = $region
$ws.Cells.Item(329, 4).Value = 44504
$ws.Cells.Item(329, 5).Value = $codreg
$ws.Cells.Item(329, 6).Value = $tipo
$ws.Cells.Item(329, 7).Value = $productoId
$ws.Cells.Item(329, 8).Value = $producto
$ws.Cells.Item(329, 9).Value = $categoriaId
$ws.Cells.Item(329, 10).Value = $categoria
$ws.Cells.Item(329, 11).Value = "Clementina"
$ws.Cells.Item(329, 12).Value = "Primera"
$ws.Cells.Item(329, 13).Value = 300
$ws.Cells.Item(329, 14).Value = 5000
$ws.Cells.Item(329, 15).Value = 5000
$ws.Cells.Item(329, 16).Value = 5000
$ws.Cells.Item(329, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(329, 18).Value = "Provincia de Petorca"
$ws.Cells.Item(329, 19).Value = 500
$ws.Cells.Item(329, 20).Value = 10

# New row 330: Clementina / Segunda
$ws.Cells.Item(330, 1).Value = $mercadoId
$ws.Cells.Item(330, 2).Value = $mercado
$ws.Cells.Item(330, 3).Value = $region
$ws.Cells.Item(330, 4).Value = 44504
$ws.Cells.Item(330, 5).Value = $codreg
$ws.Cells.Item(330, 6).Value = $tipo
$ws.Cells.Item(330, 7).Value = $productoId
$ws.Cells.Item(330, 8).Value = $producto
$ws.Cells.Item(330, 9).Value = $categoriaId
$ws.Cells.Item(330, 10).Value = $categoria
$ws.Cells.Item(330, 11).Value = "Clementina"
$ws.Cells.Item(330, 12).Value = "Segunda"
$ws.Cells.Item(330, 13).Value = 410
$ws.Cells.Item(330, 14).Value = 4000
$ws.Cells.Item(330, 15).Value = 4000
$ws.Cells.Item(330, 16).Value = 4000
$ws.Cells.Item(330, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(330, 18).Value = "Provincia de Petorca"
$ws.Cells.Item(330, 19).Value = 400
$ws.Cells.Item(330, 20).Value = 10
